$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (column I) values with NOAA data,
# and the downstream worst/best ASHP COP values (columns N and O) that
# depend on temperature, for the affected rows.

$updates = @{
    "I2"  = 19.60879629629628
    "N2"  = 1.523070892784755
    "O2"  = 1.622381208625714

    "I3"  = 19.60879629629628

    "I6"  = 17.72453703703706

    "I7"  = -0.763888888888889
    "N7"  = 1.32738293362297
    "O7"  = 1.399902772843707

    "I8"  = -0.763888888888889

    "I9"  = 15.74228395061728
    "N9"  = 0.9359630391213685
    "O9"  = 0.9581010400034263

    "I10" = 15.74228395061728

    "I11" = -3.083333333333334

    "I12" = 1.925925925925943

    "I13" = 1.925925925925943

    "I14" = 14.96875
    "N14" = 0.9337049960039765
    "O14" = 0.9557151767570068

    "I15" = 14.96875

    "I22" = 15.74228395061728
    "N22" = 0.9359630391213685
    "O22" = 0.9581010400034263

    "I23" = 15.74228395061728

    "I28" = 14.96875
    "N28" = 1.473592088566053
    "O28" = 1.565708370582976

    "I29" = 14.96875

    "I30" = 1.925925925925943

    "I36" = 14.96875
    "N36" = 0.9337049960039765
    "O36" = 0.9557151767570068

    "I37" = 14.96875
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
